# Updated symbol list with refreshed Price / Volume(1h) figures.
# Values are kept as plain text (matching the source sheet's inline-string
# cells) by forcing a text NumberFormat before assignment, so e.g. "313.51"
# and "2.65%" are not silently re-typed as numbers/percentages by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "313.51"
    "E2" = "2.65%"
    "D3" = "35.51"
    "E3" = "-2.04%"
    "E4" = "1.25%"
    "D5" = "0.08179"
    "E5" = "3.89%"
    "D6" = "2.098"
    "E6" = "-1.97%"
    "D7" = "7.976"
    "E7" = "0.45%"
    "D8" = "4.138"
    "E8" = "0.14%"
    "D9" = "0.9288"
    "E9" = "0.79%"
    "D10" = "0.1041"
    "E10" = "6.82%"
    "D11" = "0.1922"
    "E11" = "3.88%"
    "D12" = "0.09192"
    "E12" = "7.02%"
    "D13" = "0.03641"
    "E13" = "1.76%"
    "D14" = "0.09891"
    "E14" = "-0.32%"
    "D15" = "0.001447"
    "E15" = "0.42%"
    "D16" = "0.005697"
    "E16" = "0.08%"
    "D17" = "3.474"
    "E17" = "0.08%"
    "E18" = "8.11%"
    "D19" = "0.3410"
    "E19" = "1.05%"
    "D20" = "0.1302"
    "E20" = "-3.35%"
    "D21" = "5.103"
    "E21" = "-0.84%"
    "E22" = "0.20%"
    "D23" = "0.04548"
    "E23" = "-0.51%"
    "D24" = "0.001233"
    "E24" = "0.07%"
    "D25" = "0.004786"
    "E25" = "-0.45%"
    "D26" = "0.0001252"
    "E26" = "-3.67%"
    "D27" = "0.0004454"
    "E27" = "-6.26%"
    "D39" = "0.01993"
    "E39" = "7.57%"
    "D40" = "0.04910"
    "E40" = "3.93%"
    "D41" = "0.007543"
    "E41" = "-2.82%"
    "D42" = "0.1383"
    "E42" = "-0.26%"
    "D43" = "0.008122"
    "E43" = "5.26%"
    "D44" = "0.002224"
    "E44" = "2.85%"
    "E45" = "1.18%"
    "D46" = "0.00006603"
    "E46" = "3.61%"
    "E47" = "0.08%"
    "D48" = "185.91"
    "E48" = "260.25%"
    "E49" = "-10.54%"
    "D50" = "0.00002103"
    "E50" = "0.08%"
    "D51" = "0.0002003"
    "E51" = "0.08%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
